$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3410236666666666
$ws.Range("H2").Value = 1.023071
$ws.Range("I2").Value = 0.01850325494520333
$ws.Range("J2").Value = 0.01850325494520333
$ws.Range("M2").Value = 2.565830333333333
$ws.Range("N2").Value = 7.697490999999999
$ws.Range("O2").Value = 0.0934185609347503
$ws.Range("P2").Value = 0.0934185609347503
$ws.Range("Q2").Value = 0.8750088683178886
$ws.Range("R2").Value = 7.875079814860999
$ws.Range("S2").Value = 0.001728547449589697
$ws.Range("T2").Value = 0.001728547449589697
$ws.Range("G3").Value = 0.3410236666666666
$ws.Range("H3").Value = 1.023071
$ws.Range("I3").Value = 0.01850325494520333
$ws.Range("J3").Value = 0.01850325494520333
$ws.Range("O3").Value = 0.3847798091300315
$ws.Range("P3").Value = 0.3847798091300315
$ws.Range("Q3").Value = 3.604056217196555
$ws.Range("R3").Value = 32.436505954769
$ws.Range("S3").Value = 0.007119678906099647
$ws.Range("T3").Value = 0.007119678906099647
$ws.Range("G4").Value = 0.3410236666666666
$ws.Range("H4").Value = 1.023071
$ws.Range("I4").Value = 0.01850325494520333
$ws.Range("J4").Value = 0.01850325494520333
$ws.Range("M4").Value = 13.68376133333333
$ws.Range("N4").Value = 41.051284
$ws.Range("O4").Value = 0.4982080363333638
$ws.Range("P4").Value = 0.4982080363333638
$ws.Range("Q4").Value = 4.666486463684889
$ws.Range("R4").Value = 41.998378173164
$ws.Range("S4").Value = 0.009218470312025351
$ws.Range("T4").Value = 0.009218470312025351
$ws.Range("G5").Value = 0.3410236666666666
$ws.Range("H5").Value = 1.023071
$ws.Range("I5").Value = 0.01850325494520333
$ws.Range("J5").Value = 0.01850325494520333
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6480206666666667
$ws.Range("N5").Value = 1.944062
$ws.Range("O5").Value = 0.0235935936018545
$ws.Range("P5").Value = 0.0235935936018545
$ws.Range("Q5").Value = 0.2209903838224444
$ws.Range("R5").Value = 1.988913454402
$ws.Range("S5").Value = 0.0004365582774886318
$ws.Range("T5").Value = 0.0004365582774886318
$ws.Range("I6").Value = 0.2085050756621187
$ws.Range("J6").Value = 0.2085050756621187
$ws.Range("M6").Value = 2.565830333333333
$ws.Range("N6").Value = 7.697490999999999
$ws.Range("O6").Value = 0.0934185609347503
$ws.Range("P6").Value = 0.0934185609347503
$ws.Range("Q6").Value = 9.860091688405443
$ws.Range("R6").Value = 88.74082519564899
$ws.Range("S6").Value = 0.01947824411594636
$ws.Range("T6").Value = 0.01947824411594636
$ws.Range("I7").Value = 0.2085050756621187
$ws.Range("J7").Value = 0.2085050756621187
$ws.Range("O7").Value = 0.3847798091300315
$ws.Range("P7").Value = 0.3847798091300315
$ws.Range("S7").Value = 0.0802285432159128
$ws.Range("T7").Value = 0.08022854321591279
$ws.Range("I8").Value = 0.2085050756621187
$ws.Range("J8").Value = 0.2085050756621187
$ws.Range("M8").Value = 13.68376133333333
$ws.Range("N8").Value = 41.051284
$ws.Range("O8").Value = 0.4982080363333638
$ws.Range("P8").Value = 0.4982080363333638
$ws.Range("Q8").Value = 52.58459206600845
$ws.Range("R8").Value = 473.261328594076
$ws.Range("S8").Value = 0.1038789043111636
$ws.Range("T8").Value = 0.1038789043111636
$ws.Range("I9").Value = 0.2085050756621187
$ws.Range("J9").Value = 0.2085050756621187
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6480206666666667
$ws.Range("N9").Value = 1.944062
$ws.Range("O9").Value = 0.0235935936018545
$ws.Range("P9").Value = 0.0235935936018545
$ws.Range("Q9").Value = 2.490243842824222
$ws.Range("R9").Value = 22.412194585418
$ws.Range("S9").Value = 0.004919384019095951
$ws.Range("T9").Value = 0.004919384019095951
$ws.Range("G10").Value = 0.2092423333333333
$ws.Range("H10").Value = 0.627727
$ws.Range("I10").Value = 0.0113530661283407
$ws.Range("J10").Value = 0.0113530661283407
$ws.Range("M10").Value = 2.565830333333333
$ws.Range("N10").Value = 7.697490999999999
$ws.Range("O10").Value = 0.0934185609347503
$ws.Range("P10").Value = 0.0934185609347503
$ws.Range("Q10").Value = 0.5368803258841111
$ws.Range("R10").Value = 4.831922932957
$ws.Range("S10").Value = 0.001060587099906645
$ws.Range("T10").Value = 0.001060587099906645
$ws.Range("G11").Value = 0.2092423333333333
$ws.Range("H11").Value = 0.627727
$ws.Range("I11").Value = 0.0113530661283407
$ws.Range("J11").Value = 0.0113530661283407
$ws.Range("O11").Value = 0.3847798091300315
$ws.Range("P11").Value = 0.3847798091300315
$ws.Range("Q11").Value = 2.211345446261444
$ws.Range("R11").Value = 19.902109016353
$ws.Range("S11").Value = 0.00436843061790356
$ws.Range("T11").Value = 0.00436843061790356
$ws.Range("G12").Value = 0.2092423333333333
$ws.Range("H12").Value = 0.627727
$ws.Range("I12").Value = 0.0113530661283407
$ws.Range("J12").Value = 0.0113530661283407
$ws.Range("M12").Value = 13.68376133333333
$ws.Range("N12").Value = 41.051284
$ws.Range("O12").Value = 0.4982080363333638
$ws.Range("P12").Value = 0.4982080363333638
$ws.Range("Q12").Value = 2.863222150163112
$ws.Range("R12").Value = 25.768999351468
$ws.Range("S12").Value = 0.005656188782163446
$ws.Range("T12").Value = 0.005656188782163446
$ws.Range("G13").Value = 0.2092423333333333
$ws.Range("H13").Value = 0.627727
$ws.Range("I13").Value = 0.0113530661283407
$ws.Range("J13").Value = 0.0113530661283407
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6480206666666667
$ws.Range("N13").Value = 1.944062
$ws.Range("O13").Value = 0.0235935936018545
$ws.Range("P13").Value = 0.0235935936018545
$ws.Range("Q13").Value = 0.1355933563415556
$ws.Range("R13").Value = 1.220340207074
$ws.Range("S13").Value = 0.0002678596283670502
$ws.Range("T13").Value = 0.0002678596283670502
$ws.Range("G14").Value = 14.03735666666667
$ws.Range("H14").Value = 42.11207
$ws.Range("I14").Value = 0.7616386032643372
$ws.Range("J14").Value = 0.7616386032643372
$ws.Range("M14").Value = 2.565830333333333
$ws.Range("N14").Value = 7.697490999999999
$ws.Range("O14").Value = 0.0934185609347503
$ws.Range("P14").Value = 0.0934185609347503
$ws.Range("Q14").Value = 36.01747553515222
$ws.Range("R14").Value = 324.15727981637
$ws.Range("S14").Value = 0.0711511822693076
$ws.Range("T14").Value = 0.0711511822693076
$ws.Range("G15").Value = 14.03735666666667
$ws.Range("H15").Value = 42.11207
$ws.Range("I15").Value = 0.7616386032643372
$ws.Range("J15").Value = 0.7616386032643372
$ws.Range("O15").Value = 0.3847798091300315
$ws.Range("P15").Value = 0.3847798091300315
$ws.Range("Q15").Value = 148.3516468578589
$ws.Range("R15").Value = 1335.16482172073
$ws.Range("S15").Value = 0.2930631563901154
$ws.Range("T15").Value = 0.2930631563901154
$ws.Range("G16").Value = 14.03735666666667
$ws.Range("H16").Value = 42.11207
$ws.Range("I16").Value = 0.7616386032643372
$ws.Range("J16").Value = 0.7616386032643372
$ws.Range("M16").Value = 13.68376133333333
$ws.Range("N16").Value = 41.051284
$ws.Range("O16").Value = 0.4982080363333638
$ws.Range("P16").Value = 0.4982080363333638
$ws.Range("Q16").Value = 192.0838383775423
$ws.Range("R16").Value = 1728.75454539788
$ws.Range("S16").Value = 0.3794544729280114
$ws.Range("T16").Value = 0.3794544729280114
$ws.Range("G17").Value = 14.03735666666667
$ws.Range("H17").Value = 42.11207
$ws.Range("I17").Value = 0.7616386032643372
$ws.Range("J17").Value = 0.7616386032643372
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6480206666666667
$ws.Range("N17").Value = 1.944062
$ws.Range("O17").Value = 0.0235935936018545
$ws.Range("P17").Value = 0.0235935936018545
$ws.Range("Q17").Value = 9.096497225371111
$ws.Range("R17").Value = 81.86847502834
$ws.Range("S17").Value = 0.01796979167690286
$ws.Range("T17").Value = 0.01796979167690286
